$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# The "Tournament Barrier (MPI)" intro paragraph gains a new closing remark
# describing message flow on the way up/down the tree. The document's
# "_GoBack" bookmark (previously sitting at the end of the "Experimentation
# methodology" paragraph) moves into the middle of this new sentence, right
# where the live cursor was when the author last edited the text.
# ---------------------------------------------------------------------------

# 1) Append the new sentence to the end of the paragraph that currently ends
#    "...to move up and down the tree. " -- using Find/Replace (rather than
#    Range.InsertAfter) so the appended text inherits the run's existing
#    sz/szCs (12pt) run formatting instead of landing with no formatting.
$rng = $d.Content
$rng.Find.Execute( `
    "to move up and down the tree. ", `
    $false, $false, $false, $false, $false, $true, 1, $false, `
    "to move up and down the tree. While moving up the tree, the losers send messages to the winners and while moving down, the winners send messages to losers.", `
    2) | Out-Null

# 2) Mark the future run-break (the spot the new bookmark will sit at) with a
#    throwaway bookmark so inserting "_GoBack" there splits the run cleanly
#    instead of merging "While moving up the tree..." back into the
#    preceding run.
$splitRng = $d.Content
$splitRng.Find.Execute( `
    "While moving up the tree", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$splitRng.Collapse(1) # wdCollapseStart
$d.Bookmarks.Add("ZZTmpSplit", $splitRng) | Out-Null

# 3) Re-seat "_GoBack" (removing it from wherever it was before -- bookmark
#    names are unique, so Bookmarks.Add on an existing name relocates it) to
#    sit between the two new sentences.
$goBackRng = $d.Content
$goBackRng.Find.Execute( `
    "while moving down,", `
    $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$goBackRng.Collapse(0) # wdCollapseEnd
$d.Bookmarks.Add("_GoBack", $goBackRng) | Out-Null

# 4) Drop the throwaway helper bookmark; the run split it created stays put.
$d.Bookmarks.Item("ZZTmpSplit").Delete()
